$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(84).Insert()

$ws.Range("A84").Value = 4
$ws.Range("B84").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C84").Value = "Los Lagos"
$ws.Range("D84").Value = 44483
$ws.Range("E84").Value = 10
$ws.Range("F84").Value = 100114014
$ws.Range("G84").Value = "Betarraga"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 1000
$ws.Range("L84").Value = 1000
$ws.Range("M84").Value = 1000
$ws.Range("N84").Value = "$/paquete 5 unidades"
$ws.Range("O84").Value = "Región del Maule"
$ws.Range("P84").Value = 200
$ws.Range("Q84").Value = 5
$ws.Range("R84").Value = "Hortaliza"
